# Removed Test Case Inter-Dependency
# - productname (B1 on both sheets) gets a "-1st" suffix so it no longer
#   collides with other automated test-case runs.
# - shortname (B2 on ProductLoanInput) switches from the numeric 2445 to
#   the text token "244d" for the same reason.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$wsInput.Range("B1").Value = "2445-RBI-EPP-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-INST-1st"
$wsInput.Range("B2").Value = "244d"

$wsOutput.Range("B1").Value = "2445-RBI-EPP-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-INST-1st"

# Restore the cursor/selection to B1 on both sheets (matches the author's
# commit), re-activating ProductLoanInput last so it stays the selected tab.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()

$wsInput.Activate()
$wsInput.Range("B1").Select()
